$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: fill in the previously-empty Status cell ---
$ws.Range("C2").Value = "Done"

# --- Row 3: brand-new task row ---
# Copy the date formatting (number format + center/center alignment) from A2
# onto A3 before writing the date value, so the new date cell matches the
# existing "Date" column styling exactly.
$ws.Range("A2").Copy()
$ws.Range("A3").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("A3").Value = 45607        # 11/11/2024
$ws.Range("B3").Value = "Making csv_operations module."
$ws.Range("C3").Value = "In progress"
$ws.Range("D3").Value = "Drew Hutchinson"

# New data made columns A (Date) and C (Status) a bit wider; re-fit just
# those two columns (B and D already accommodate the new content).
$ws.Columns("A").ColumnWidth = 10.66
$ws.Columns("C").ColumnWidth = 9.96

# The user ended up with the cursor on B4 after entering the new row.
[void]$ws.Range("B4").Select()
